$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Code Your Own Linear Regression" is being added as new in-class exercise
# content for the Parallelism class session (row 24). The exercise-link
# string needs to be written before the "Do Before Class" text so that the
# shared-string table ends up with the two new strings in the same order as
# the target workbook (link string, then the updated bullet-list string).
$ws.Range("D24").Value = "``Link <exercises/Exercise_codeyourownlinearregression.ipynb>``_"
$ws.Range("C24").Value = "- ``Parallel Computing <parallelism.ipynb>``_`n- Review Linear Regression Matrix Algebra`n- Review Defining Classes"

# Row grew from one line to three lines of wrapped text, so it needs to be
# taller to show everything.
$ws.Rows.Item(24).RowHeight = 51

# Reflect where the cursor ended up after making the edit.
$ws.Range("C25").Select()
